# Update the "Förändrad" (Changed) date column (C) for rows 2-5
# from serial date 45243 (2023-11-13) to 45244 (2023-11-14).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 45244
$ws.Range("C3").Value = 45244
$ws.Range("C4").Value = 45244
$ws.Range("C5").Value = 45244
